# Refresh the crypto price/volume table in place (values only).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.242.01'
$ws.Range("E2").Value = '  -0.42%  '

$ws.Range("D3").Value = '1.588.88'
$ws.Range("E3").Value = '  -0.22%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").Value = '''211.84'
$ws.Range("E5").Value = '  +0.80%  '

$ws.Range("D6").Value = '''0.503'
$ws.Range("E6").Value = '  +0.16%  '

$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").Value = '''0.246'
$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").Value = '  -0.96%  '

$ws.Range("D10").Value = '''19.25'
$ws.Range("E10").Value = '  -1.92%  '

$ws.Range("E11").Value = '  +0.41%  '

$ws.Range("D12").Value = '1.811.35'
$ws.Range("E12").Value = '  -0.21%  '

$ws.Range("D13").Value = '1.583.39'
$ws.Range("E13").Value = '  -0.75%  '

$ws.Range("E14").Value = '  -1.46%  '

$ws.Range("E15").Value = '  -0.18%  '

$ws.Range("E16").Value = '  -1.00%  '

$ws.Range("D17").Value = '26.226.39'
$ws.Range("E17").Value = '  -0.47%  '

$ws.Range("E18").Value = '  -0.71%  '

$ws.Range("D19").Value = '''215.26'
$ws.Range("E19").Value = '  +1.37%  '

$ws.Range("E20").Value = '  -1.31%  '

$ws.Range("D21").Value = '''0.999'
$ws.Range("E21").Value = '  -0.06%  '

$ws.Range("D22").Value = '''4.25'
$ws.Range("E22").Value = '  -1.01%  '

$ws.Range("D23").Value = '''2.17'
$ws.Range("E23").Value = '  -0.63%  '

$ws.Range("E24").Value = '  +0.48%  '

$ws.Range("D25").Value = '''144.04'
$ws.Range("E25").Value = '  -0.53%  '

$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.09%  '

$ws.Range("D27").Value = '''6.99'
$ws.Range("E27").Value = '  -0.92%  '

$ws.Range("E28").Value = '  -0.58%  '

$ws.Range("E29").Value = '  -0.93%  '

$ws.Range("E30").Value = '  -1.58%  '

$ws.Range("E31").Value = '  +0.55%  '

$ws.Range("E32").Value = '  -0.95%  '

$ws.Range("D33").Value = '1.365.34'
$ws.Range("E33").Value = '  +5.22%  '

$ws.Range("E34").Value = '  -1.98%  '

$ws.Range("E35").Value = '  -0.38%  '

$ws.Range("E36").Value = '  -1.27%  '

$ws.Range("E37").Value = '  -4.95%  '

$ws.Range("E38").Value = '  -0.68%  '

$ws.Range("E39").Value = '  +0.71%  '

$ws.Range("E40").Value = '  +3.64%  '

$ws.Range("E41").Value = '  -0.15%  '

$ws.Range("E42").Value = '  +0.87%  '

$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").Value = '''2.14'
$ws.Range("E43").Value = '  -0.07%  '

$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").Value = '''0.923'
$ws.Range("E44").Value = '  -17.20%  '

$ws.Range("D45").Value = '1.723.83'
$ws.Range("E45").Value = '  -0.19%  '

$ws.Range("D46").Value = '''61.02'
$ws.Range("E46").Value = '  -2.88%  '

$ws.Range("D47").Value = '''86.22'
$ws.Range("E47").Value = '  -2.56%  '

$ws.Range("E48").Value = '  -1.10%  '

$ws.Range("E49").Value = '  -1.94%  '

$ws.Range("D50").Value = '''0.0984'
$ws.Range("E50").Value = '  -1.81%  '

$ws.Range("E51").Value = '  -1.17%  '
